$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template row used purely as a style donor (its own cells/hyperlinks are left untouched).
$styleDonorRow = 487

function Set-EventRow($Row, $Date, $Event, $Location, $City, $Link) {

    $bcdSrc = "B$($styleDonorRow):D$($styleDonorRow)"
    $bcdDst = "B$($Row):D$($Row)"
    $eSrc = "E$($styleDonorRow)"
    $eDst = "E$($Row)"

    # Pull in the plain "text" cell style (s=3) used by every populated row for B:D.
    $ws.Range($bcdSrc).Copy($ws.Range($bcdDst))

    $ws.Range("A$($Row)").Value = $Date
    $ws.Range("B$($Row)").Value = $Event
    $ws.Range("C$($Row)").Value = $Location
    $ws.Range("D$($Row)").Value = $City

    # Register the external hyperlink relationship for column E first (the share-string
    # entry for the plain URL text it allocates gets overwritten a few lines down once the
    # B:D text above has already claimed the earlier shared-string slots).
    $ws.Hyperlinks.Add($ws.Range($eDst), $Link, "", "", $Link) | Out-Null

    # Pull in the plain "text" cell style (s=3) for E too (Hyperlinks.Add switches it to a
    # dedicated hyperlink style), then write the URL text and re-apply the blue/underline
    # run formatting character-by-character so it lands as rich text inside the shared
    # string (matching every other link cell in the sheet) instead of a cell-level style.
    $ws.Range($eSrc).Copy($ws.Range($eDst))
    $ws.Range($eDst).Value = $Link
    $len = $Link.Length
    $ws.Range($eDst).Characters(1, $len - 1).Font.Underline = 2
    $ws.Range($eDst).Characters(1, $len - 1).Font.Color = 65280
    $ws.Range($eDst).Characters($len, 1).Font.Underline = 2
    $ws.Range($eDst).Characters($len, 1).Font.Color = 65280
}

Set-EventRow 488 45850 "BLACK SECTOR" "Elektroküche" "Köln" "https://www.instagram.com/reel/DKFEgC4sMyS/?igsh=Z3lpNmdla2ozbGI3"
Set-EventRow 489 45822 "NOCTURGENERATION x OUTLAW RAVES" "check event link" "Köln" "https://www.instagram.com/reel/DJrrBO9sfDQ/?igsh=MWV3YWI3ZmRhcDNucA=="
Set-EventRow 490 45813 "240 MIN RAVE" "SNRS" "Dortmund" "https://www.instagram.com/p/DKH2hI7IKZt/?igsh=bmhjMG51d3ZrZzh6"
Set-EventRow 491 45805 "SBA" "Mikroport" "Krefeld" "https://www.instagram.com/p/DKFDZizI3wk/?igsh=MTVwbmEwY3B1YnluYg=="
Set-EventRow 492 45857 "FASTER DAY & NIGHT w/ KLANGKUENSTLER" "Junkyard" "Dortmund" "https://www.instagram.com/reel/DKFRwmRqJeB/?igsh=MWlrN3Q1dWRzbzVzOA=="
Set-EventRow 493 45842 "NEON DREAMS" "Artheater" "Köln" "https://www.instagram.com/p/DKFSde_svqT/?igsh=MTd6YXFpOWd3Zno1ag=="
Set-EventRow 494 45961 "UNREAL WAREHOUSE" "Jahrhunderthalle" "Bochum" "https://www.instagram.com/unrealgermany?igsh=MWgxbHJhYzAyMHI0Zw=="

Write-Output "done"
